# "Add files via upload" — SQL-Data.xlsx
#
# Net effect per the diff:
#  - Ch05-Assignment!E13:E32: the hire-date column switches from real
#    date/time serials (numFmtId 22, "m/d/yyyy h:mm") to literal text
#    strings (numFmtId 49, "@"), one new shared string per row.
#  - A stray/unused extra column (originally column K, with no data in
#    it) is removed from Ch05-Assignment, which re-numbers the bestFit
#    width hints that used to sit at columns L/M down to K/L.
#  - The previously-active sheet/tab moves from "Chapter 05" to
#    "Ch05-Assignment", with a new selected cell on that sheet.

$wb = $excel.ActiveWorkbook

$wsAssignment = $wb.Worksheets.Item("Ch05-Assignment")

# --- Remove the empty, data-less column that used to sit at column K.
# There is no cell data out past column J on this sheet (dimension is
# A1:J32), so this only renumbers the leftover bestFit column-width
# hints that lived at columns L/M down to columns K/L.
$wsAssignment.Columns.Item(11).Delete()

# --- Replace the numeric hire-date values in E13:E32 with literal text
# timestamps (including the one data-entry glitch in the source row
# that has a trailing ",60000" baked into the text).
$hireDateText = @(
    "2020-05-12 09:00:00",
    "2021-08-22 14:30:00",
    "2022-09-10 10:15:00",
    "2023-03-18 13:45:00",
    "2020-12-01 08:30:00",
    "2021-06-25 11:00:00",
    "2022-10-05 16:00:00",
    "2023-04-12 12:30:00",
    "2021-09-15 09:45:00",
    "2022-11-20 10:20:00",
    "2023-05-10 09:30:00",
    "2022-08-15 14:00:00",
    "2021-12-01 10:45:00",
    "2020-07-20 13:20:00",
    "2022-03-18 11:10:00",
    "2023-01-25 15:00:00",
    "2021-10-30 12:00:00,60000",
    "2022-06-12 09:00:00",
    "2023-02-14 16:30:00",
    "2020-09-05 08:15:00"
)

$row = 13
foreach ($txt in $hireDateText) {
    $cell = $wsAssignment.Cells.Item($row, 5)
    $cell.Value = $txt
    $cell.NumberFormat = "@"
    $row++
}

# --- Move the active tab / selection from "Chapter 05" to
# "Ch05-Assignment" (matches the workbook.xml activeTab flip and the
# tabSelected/selection swap between the two sheetViews).
$wsAssignment.Activate()
$wsAssignment.Range("L10").Select()
